{"js": "// Replace the date line and all of the \"NNN\u00d7N=\" problem cells with their\n// updated values. Every old value in this document is unique, so a plain\n// body-wide search-and-replace (matchCase, whole match) for each pair is\n// unambiguous and order-independent.\nconst replacements = [\n  [\"2024-06-08 Saturday\", \"2024-06-09 Sunday\"],\n  [\"914\u00d79=\", \"898\u00d78=\"],\n  [\"509\u00d77=\", \"611\u00d73=\"],\n  [\"281\u00d74=\", \"338\u00d78=\"],\n  [\"822\u00d72=\", \"704\u00d73=\"],\n  [\"338\u00d77=\", \"823\u00d75=\"],\n  [\"112\u00d74=\", \"513\u00d75=\"],\n  [\"352\u00d77=\", \"652\u00d77=\"],\n  [\"207\u00d75=\", \"196\u00d75=\"],\n  [\"631\u00d77=\", \"712\u00d78=\"],\n  [\"696\u00d72=\", \"995\u00d79=\"],\n  [\"840\u00d72=\", \"881\u00d75=\"],\n  [\"133\u00d77=\", \"912\u00d76=\"],\n  [\"641\u00d78=\", \"175\u00d76=\"],\n  [\"131\u00d74=\", \"514\u00d72=\"],\n  [\"385\u00d72=\", \"184\u00d72=\"],\n  [\"990\u00d74=\", \"770\u00d73=\"],\n  [\"777\u00d79=\", \"817\u00d77=\"],\n  [\"567\u00d74=\", \"349\u00d78=\"],\n  [\"755\u00d75=\", \"180\u00d72=\"],\n  [\"675\u00d79=\", \"245\u00d73=\"],\n  [\"730\u00d73=\", \"904\u00d74=\"],\n  [\"472\u00d72=\", \"427\u00d78=\"],\n  [\"961\u00d75=\", \"569\u00d79=\"],\n  [\"124\u00d77=\", \"785\u00d79=\"],\n  [\"686\u00d76=\", \"454\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and all of the \"NNN\u00d7N=\" problem cells with their\n# updated values. Every \"find\" value in this document is unique, so a\n# whole-document Find/Replace (MatchCase on, MatchWholeWord off \u2014 the text\n# itself is the whole run) for each pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-06-08 Saturday\", \"2024-06-09 Sunday\"),\n    @(\"914\u00d79=\", \"898\u00d78=\"),\n    @(\"509\u00d77=\", \"611\u00d73=\"),\n    @(\"281\u00d74=\", \"338\u00d78=\"),\n    @(\"822\u00d72=\", \"704\u00d73=\"),\n    @(\"338\u00d77=\", \"823\u00d75=\"),\n    @(\"112\u00d74=\", \"513\u00d75=\"),\n    @(\"352\u00d77=\", \"652\u00d77=\"),\n    @(\"207\u00d75=\", \"196\u00d75=\"),\n    @(\"631\u00d77=\", \"712\u00d78=\"),\n    @(\"696\u00d72=\", \"995\u00d79=\"),\n    @(\"840\u00d72=\", \"881\u00d75=\"),\n    @(\"133\u00d77=\", \"912\u00d76=\"),\n    @(\"641\u00d78=\", \"175\u00d76=\"),\n    @(\"131\u00d74=\", \"514\u00d72=\"),\n    @(\"385\u00d72=\", \"184\u00d72=\"),\n    @(\"990\u00d74=\", \"770\u00d73=\"),\n    @(\"777\u00d79=\", \"817\u00d77=\"),\n    @(\"567\u00d74=\", \"349\u00d78=\"),\n    @(\"755\u00d75=\", \"180\u00d72=\"),\n    @(\"675\u00d79=\", \"245\u00d73=\"),\n    @(\"730\u00d73=\", \"904\u00d74=\"),\n    @(\"472\u00d72=\", \"427\u00d78=\"),\n    @(\"961\u00d75=\", \"569\u00d79=\"),\n    @(\"124\u00d77=\", \"785\u00d79=\"),\n    @(\"686\u00d76=\", \"454\u00d73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
